$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 99426637
$ws.Range("B2").Value = 78569
$ws.Range("E2").Value = 6458
$ws.Range("F2").Value = "Lunglav"
$ws.Range("G2").Value = "Lobaria pulmonaria"
$ws.Range("H2").Value = "(L.) Hoffm."
$ws.Range("M2").Value = ""
$ws.Range("Q2").Value = 626058.8445098634
$ws.Range("R2").Value = 7010453.642544731
$ws.Range("S2").Value = 25

$ws.Range("A3").Value = 99426643
$ws.Range("B3").Value = 78527
$ws.Range("D3").Value = "LC"
$ws.Range("E3").Value = 229497
$ws.Range("F3").Value = "Korallblylav"
$ws.Range("G3").Value = "Parmeliella triptophylla"
$ws.Range("H3").Value = "(Ach.) Müll.Arg."
$ws.Range("Q3").Value = 625950.5816682897
$ws.Range("R3").Value = 7010469.289106234

$ws.Range("A4").Value = 99426638
$ws.Range("B4").Value = 78569
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 6458
$ws.Range("F4").Value = "Lunglav"
$ws.Range("G4").Value = "Lobaria pulmonaria"
$ws.Range("H4").Value = "(L.) Hoffm."
$ws.Range("Q4").Value = 625948.7751551091
$ws.Range("R4").Value = 7010469.218583253

$ws.Range("A5").Value = 99426616
$ws.Range("B5").Value = 56395
$ws.Range("E5").Value = 100109
$ws.Range("F5").Value = "Tretåig hackspett"
$ws.Range("G5").Value = "Picoides tridactylus"
$ws.Range("H5").Value = "(Linnaeus, 1758)"
$ws.Range("M5").Value = "äldre spår"
$ws.Range("Q5").Value = 626303.9086833318
$ws.Range("R5").Value = 7010513.817294765

$ws.Range("A6").Value = 99426608
$ws.Range("B6").Value = 89392
$ws.Range("E6").Value = 1202
$ws.Range("F6").Value = "Ullticka"
$ws.Range("G6").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H6").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("M6").Value = ""
$ws.Range("Q6").Value = 625848.4141095353
$ws.Range("R6").Value = 7010491.0488934

$ws.Range("A7").Value = 99426622
$ws.Range("B7").Value = 56411
$ws.Range("E7").Value = 100049
$ws.Range("F7").Value = "Spillkråka"
$ws.Range("G7").Value = "Dryocopus martius"
$ws.Range("H7").Value = "(Linnaeus, 1758)"
$ws.Range("M7").Value = "spel/sång"
$ws.Range("Q7").Value = 626015.2164001387
$ws.Range("R7").Value = 7010979.93372071
$ws.Range("S7").Value = 98
